# "May Month code updated" -- roll every report's From/End Date window
# forward by one month, fix a couple of Duration labels, and add a new
# "Unit wise Plan Qty" report row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Buyer wise monthly plan qty ---------------------------------
$ws.Range("C2").Value = 45809
$ws.Range("D2").Value = 45961

# --- Row 3: Group Booking ------------------------------------------------
$ws.Range("C3").Value = 45778
$ws.Range("D3").Value = 45961

# --- Row 4: Provision ------------------------------------------------
$ws.Range("C4").Value = 45778
$ws.Range("D4").Value = 45961
$ws.Range("E4").Value = "Current Month+ 5"

# --- Row 5: Monthly Blank Days -------------------------------------------
$ws.Range("C5").Value = 45809
$ws.Range("D5").Value = 45961

# --- Row 6: unit+ buyer wise report --------------------------------------
$ws.Range("C6").Value = 45809
$ws.Range("D6").Value = 45961

# --- Row 7: Factory wise Weekly blank days -------------------------------
$ws.Range("C7").Value = 45809
$ws.Range("D7").Value = 45900

# --- Row 8: Plan and Efficiency Report -----------------------------------
$ws.Range("C8").Value = 45809
$ws.Range("D8").Value = 45838

# --- Row 9: Per_Day_Requirement (Unit) -----------------------------------
$ws.Range("C9").Value = 45809
$ws.Range("D9").Value = 45838

# --- Row 10: BuyerWise Monthly Print Req. --------------------------------
$ws.Range("C10").Value = 45809
$ws.Range("D10").Value = 45900

# --- Row 11: BuyerWise Monthly Emb Req -----------------------------------
$ws.Range("C11").Value = 45809
$ws.Range("D11").Value = 45900

# --- Row 12: BuyerWise Monthly Wash Req ----------------------------------
$ws.Range("C12").Value = 45809
$ws.Range("D12").Value = 45900

# --- Row 13: Full Plan -- duration label changes from "Current Month"
#     to "Next Month", along with the date shift --------------------------
$ws.Range("C13").Value = 45809
$ws.Range("D13").Value = 45838
$ws.Range("E13").Value = "Next Month"

# --- Row 14 (new): Unit wise Plan Qty ------------------------------------
# Clone row 13's formatting (styles + row height) onto the new row, then
# fill in its values.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(14).RowHeight = 30

$ws.Range("A14").Value = "Unit wise Plan Qty"
$ws.Range("B14").Value = "Reports-Planning-Buyer wise plan qty-Unit wise plan qty"
$ws.Range("C14").Value = 45809
$ws.Range("D14").Value = 45838
$ws.Range("E14").Value = "Next Month"

# --- View state: selection moves to D2, scroll back to the top-left -----
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
